$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 71: Course, Hours, Notes for the JS101 multiply-lists / list-of-digits entry
$ws.Range("B71").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C71").Value = 0.5
$ws.Range("D71").Value = "2 small problems"

# Update the selected cell shown in the sheet view
$ws.Range("E71").Select()

$wb.Save()
